$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$brk = [char]11

$t.Cell(1,1).Range.Text = "33 x 47" + $brk + "  4    7" + $brk + "  ----" + $brk + "3|    |" + $brk + "3|    |"
$t.Cell(1,2).Range.Text = "59 x 50" + $brk + "  5    0" + $brk + "  ----" + $brk + "5|    |" + $brk + "9|    |"
$t.Cell(1,3).Range.Text = "71 x 84" + $brk + "  8    4" + $brk + "  ----" + $brk + "7|    |" + $brk + "1|    |"

$t.Cell(2,1).Range.Text = "42 x 51" + $brk + "  5    1" + $brk + "  ----" + $brk + "4|    |" + $brk + "2|    |"
$t.Cell(2,2).Range.Text = "68 x 41" + $brk + "  4    1" + $brk + "  ----" + $brk + "6|    |" + $brk + "8|    |"
$t.Cell(2,3).Range.Text = "18 x 11" + $brk + "  1    1" + $brk + "  ----" + $brk + "1|    |" + $brk + "8|    |"

$t.Cell(3,1).Range.Text = "17 x 54" + $brk + "  5    4" + $brk + "  ----" + $brk + "1|    |" + $brk + "7|    |"
$t.Cell(3,2).Range.Text = "12 x 12" + $brk + "  1    2" + $brk + "  ----" + $brk + "1|    |" + $brk + "2|    |"
$t.Cell(3,3).Range.Text = "91 x 59" + $brk + "  5    9" + $brk + "  ----" + $brk + "9|    |" + $brk + "1|    |"

$t.Cell(4,1).Range.Text = "98 x 47" + $brk + "  4    7" + $brk + "  ----" + $brk + "9|    |" + $brk + "8|    |"
$t.Cell(4,2).Range.Text = "57 x 79" + $brk + "  7    9" + $brk + "  ----" + $brk + "5|    |" + $brk + "7|    |"
$t.Cell(4,3).Range.Text = "78 x 12" + $brk + "  1    2" + $brk + "  ----" + $brk + "7|    |" + $brk + "8|    |"

$t.Cell(5,1).Range.Text = "18 x 91" + $brk + "  9    1" + $brk + "  ----" + $brk + "1|    |" + $brk + "8|    |"
$t.Cell(5,2).Range.Text = "88 x 64" + $brk + "  6    4" + $brk + "  ----" + $brk + "8|    |" + $brk + "8|    |"
$t.Cell(5,3).Range.Text = "73 x 83" + $brk + "  8    3" + $brk + "  ----" + $brk + "7|    |" + $brk + "3|    |"

Write-Output "done"
